# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" column (E16:E60) previously listed 45 period codes
# in descending order (2003 .. 1607). The update lists them in ascending
# order (1607 .. 2003) and the matching "Salario Basico" values (F16:F60)
# are reassigned so that the first 26 rows (1607..1808) carry 27580 and
# the remaining 19 rows (1809..2003) carry 31249.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periods = @(
    "1607","1608","1609","1610","1611","1612",
    "1701","1702","1703","1704","1705","1706","1707","1708","1709","1710","1711","1712",
    "1801","1802","1803","1804","1805","1806","1807","1808","1809","1810","1811","1812",
    "1901","1902","1903","1904","1905","1906","1907","1908","1909","1910","1911","1912",
    "2001","2002","2003"
)

$firstRow = 16

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = $firstRow + $i

    # Column E: "Periodo Mora" - must stay text (e.g. "1607"), not a number.
    $ws.Cells.Item($row, 5).Value = $periods[$i]

    # Column F: "Salario Basico"
    if ($row -le 41) {
        $ws.Cells.Item($row, 6).Value = 27580
    } else {
        $ws.Cells.Item($row, 6).Value = 31249
    }
}
